# Bigger banner indicators on PDF report
#
# The presentation gets a brand new "banner" slide inserted right after the
# current slide 2. The new slide is a near-duplicate of slide 2's banner
# group (blue rectangle + white rectangle), just enlarged. We build it by
# duplicating slide 2 (so the new slide lands immediately after it, with the
# correct new SlideID) and then resizing/renaming its copied shapes to match
# the bigger banner dimensions.

$p = $ppt.ActivePresentation

# Slide 2 holds the template banner group ("Group 1": Rectangle 3 + Rectangle 2).
$sourceSlide = $p.Slides.Item(2)

# Duplicating inserts the new slide immediately after the source slide, i.e.
# at index 3 -- exactly where the new banner slide belongs.
$newSlide = $sourceSlide.Duplicate().Item(1)

# Grab the duplicated banner group on the new slide.
$bannerGroup = $newSlide.Shapes.Item(1)

# Ungroup so we can resize the individual rectangles, then rebuild the group
# so the resulting group's child-extent matches its outer extent 1:1 (no
# residual scale factor), same as the original template group.
$parts = $bannerGroup.Ungroup()

$rect3 = $newSlide.Shapes.Item(1)
$rect2 = $newSlide.Shapes.Item(2)

# Rectangle 3 (the big blue rectangle) keeps its original size/position.
# Rectangle 2 (the white card) grows to the new, bigger banner size.
$rect2.Left = 337.7431496062992
$rect2.Top = 0
$rect2.Width = 374.47661417322837
$rect2.Height = 201.00141732283464

$bannerRange = $newSlide.Shapes.Range(@($rect3.Name, $rect2.Name))
$newGroup = $bannerRange.Group()
$newGroup.Name = "Group 4"
